# CIERRE 8 NOV 2021
# Fill in the previously-blank credit rows (35-39) with their actual
# remision dates / client / importe data, flag row 38 with the
# "falta la rosa" note (highlighted in a new purple fill), and move the
# active selection to B42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 35 : OBRADOR, 2515 ----
$ws.Range("A35").Value = 44498
$ws.Range("D35").Value = "OBRADOR"
$ws.Range("E35").Value = 2515

# ---- Row 36 : OBRADOR, 340 ----
$ws.Range("A36").Value = 44499
$ws.Range("D36").Value = "OBRADOR"
$ws.Range("E36").Value = 340

# ---- Row 37 : COMERCIO CENTRAL, 15657 ----
$ws.Range("A37").Value = 44499
$ws.Range("D37").Value = "COMERCIO   CENTRAL "
$ws.Range("E37").Value = 15657

# ---- Row 38 : CANCELADA, 0, flagged "falta la rosa" ----
$ws.Range("A38").Value = 44500
$ws.Range("D38").Value = "CANCELADA"
$ws.Range("D38").Font.Bold = $true
$ws.Range("D38").Font.Size = 12
$ws.Range("D38").Font.Color = 255
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = "falta la rosa"
$ws.Range("F38").Interior.Color = 16737945
$ws.Range("G38").Interior.Color = 16737945

# ---- Row 39 : COMERCIO CENTRAL, 8585 ----
$ws.Range("A39").Value = 44500
$ws.Range("D39").Value = "COMERCIO   CENTRAL "
$ws.Range("E39").Value = 8585

# ---- move the active cell/selection to B42 ----
$ws.Range("B42").Select()
